# "fix unit test done"
# The "../sample/Lib/test " placeholder-path text on the folded-corner
# shapes of slides 2 and 3 had a trailing stray space and was missing the
# "er" suffix ("test " -> "tester"). Fix both occurrences.

$p = $ppt.ActivePresentation

function Get-LibPathShape($slide) {
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t -like "*sample/Lib/test*") {
                return $shp
            }
        }
    }
    return $null
}

# Slide 2 ("S101"): text is split as "../sample/Lib" + "/tester" across two
# runs (the second run carries no explicit proofing-dirty flag).
$s2 = $p.Slides.Item(2)
$shp2 = Get-LibPathShape $s2
$tr2 = $shp2.TextFrame.TextRange
$fix2 = $tr2.Characters(14, 6)
$fix2.Text = "/tester"

# Slide 3 ("S102"): simple in-place fix, single run.
$s3 = $p.Slides.Item(3)
$shp3 = Get-LibPathShape $s3
$tr3 = $shp3.TextFrame.TextRange
$fix3 = $tr3.Characters(1, 19)
$fix3.Text = "../sample/Lib/tester"
